$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.289.61'
$ws.Range("E2").Value = '  +0.26%  '

$ws.Range("D3").Value = '3.372.98'
$ws.Range("E3").Value = '  +1.83%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.34'
$ws.Range("E5").Value = '  +0.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.61'
$ws.Range("E6").Value = '  +8.81%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.372.85'
$ws.Range("E8").Value = '  +1.87%  '

$ws.Range("E9").Value = '  -0.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.59'
$ws.Range("E10").Value = '  +5.95%  '

$ws.Range("E11").Value = '  +4.92%  '

$ws.Range("E12").Value = '  +4.86%  '

$ws.Range("D13").Value = '3.937.97'
$ws.Range("E13").Value = '  +1.45%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.121'
$ws.Range("E14").Value = '  +2.42%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000172'
$ws.Range("E15").Value = '  +2.68%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.364.38'
$ws.Range("E16").Value = '  +1.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.25'
$ws.Range("E17").Value = '  +2.21%  '

$ws.Range("D18").Value = '61.300.45'
$ws.Range("E18").Value = '  +0.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.99'
$ws.Range("E19").Value = '  +6.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.80'
$ws.Range("E20").Value = '  +4.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.43'
$ws.Range("E21").Value = '  +4.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '382.69'
$ws.Range("E22").Value = '  +8.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.575'
$ws.Range("E23").Value = '  +3.97%  '

$ws.Range("D24").Value = '3.505.18'
$ws.Range("E24").Value = '  +1.76%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.85'
$ws.Range("E26").Value = '  +1.02%  '

$ws.Range("E27").Value = '  +10.94%  '

$ws.Range("E28").Value = '  +12.71%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.79'
$ws.Range("E29").Value = '  +9.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.03%  '

$ws.Range("E31").Value = '  +3.25%  '

$ws.Range("E32").Value = '  +5.63%  '

$ws.Range("E33").Value = '  +1.20%  '

$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("D35").Value = '3.398.54'
$ws.Range("E35").Value = '  +1.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.46'
$ws.Range("E36").Value = '  +5.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.55'
$ws.Range("E37").Value = '  +2.88%  '

$ws.Range("E38").Value = '  +3.61%  '

$ws.Range("E39").Value = '  +4.34%  '

$ws.Range("E40").Value = '  +0.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0802'
$ws.Range("E41").Value = '  +6.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.41'
$ws.Range("E43").Value = '  +4.24%  '

$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.53'
$ws.Range("E44").Value = '  +1.62%  '

$ws.Range("E45").Value = '  +1.89%  '

$ws.Range("E46").Value = '  +8.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.63'
$ws.Range("E47").Value = '  +5.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.25'
$ws.Range("E48").Value = '  +3.70%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.97'
$ws.Range("E49").Value = '  +5.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.21'
$ws.Range("E50").Value = '  +12.37%  '

$ws.Range("E51").Value = '  +11.74%  '
